# Logged Week 15 and simulated Week 16
# Update the "R" (row 3) target-depth totals on both the OFF and DEF sheets
# to reflect the newly logged/simulated week's cumulative numbers.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 507
$wsOff.Range("C3").Value = 359
$wsOff.Range("D3").Value = 131

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 351
$wsDef.Range("C3").Value = 225
$wsDef.Range("D3").Value = 89
$wsDef.Range("E3").Value = 37
$wsDef.Range("F3").Value = 9
$wsDef.Range("G3").Value = 6
